$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: extend header row with Aug..Dec (columns I:M) ---
$ws1.Range("I1").Value = "Aug"
$ws1.Range("J1").Value = "Sep"
$ws1.Range("K1").Value = "Oct"
$ws1.Range("L1").Value = "Nov"
$ws1.Range("M1").Value = "Dec"

# --- Sheet1: row 2 (hansraj) gets monthly attendance numbers, B2 reset to 0 ---
$ws1.Range("B2").Value = 0
$ws1.Range("I2").Value = 5
$ws1.Range("J2").Value = 13
$ws1.Range("K2").Value = 12
$ws1.Range("L2").Value = 2
$ws1.Range("M2").Value = 6

# --- Sheet1: rows 3-7, column B reset to 0 (Jan attendance cleared for others) ---
$ws1.Range("B3").Value = 0
$ws1.Range("B4").Value = 0
$ws1.Range("B5").Value = 0
$ws1.Range("B6").Value = 0
$ws1.Range("B7").Value = 0
